$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Append the "(This is a change - Version for branch alternate)"
#    annotation to the end of the first paragraph.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.First
$insertPos = $p1.Range.End - 1          # just before the paragraph mark

# two plain spaces, no special formatting (matches the surrounding run)
$spacer = $d.Range($insertPos, $insertPos)
$spacer.InsertAfter("  ")

# the red annotation text, split the way it was originally typed/pasted
$dash = [char]0x2013
$part1 = "(This is a change " + $dash + " Ve"
$part2 = "rsion for branch alternate"
$part3 = ")"

$r1 = $d.Range($spacer.End, $spacer.End)
$r1.InsertAfter($part1)
$r1.Font.Color = 192          # RGB(192,0,0) -> w:color val="C00000"

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter($part2)
$r2.Font.Color = 192

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter($part3)
$r3.Font.Color = 192

# ---------------------------------------------------------------------
# 2) Add a new, completely empty paragraph at the very end of the body
#    (right before the sectPr), i.e. a bare <w:p/>.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$blankParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$newPara.Range.InsertXML($blankParagraphXml)

Write-Host "Done."
